$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs(1).Range.Text = "2025-11-19 Wednesday"

# New values for the 20x5 practice table, in row-major reading order
$newValues = @(
    "76-37=", "86-78=", "61-53=", "24-15=", "28+45=",
    "98-26=", "29+55=", "70+28=", "69+19=", "61-16=",
    "71-26=", "87-42=", "30+8=", "34+48=", "65-12=",
    "22-18=", "97-28=", "39-5=", "23+6=", "29-4=",
    "37+55=", "84-20=", "76-10=", "49-26=", "73+3=",
    "83-5=", "16+39=", "67-11=", "18+55=", "25-17=",
    "89-76=", "6-1=", "70+20=", "60-7=", "18+33=",
    "77-72=", "55-54=", "71-20=", "85-23=", "15+67=",
    "11+49=", "39+50=", "91-0=", "93-46=", "4+53=",
    "52-23=", "41-3=", "94-45=", "36-14=", "27+21=",
    "87-22=", "11+43=", "11+31=", "46-33=", "32+22=",
    "54-47=", "25+9=", "43-19=", "64-3=", "3+20=",
    "28-3=", "48+35=", "52+10=", "81-32=", "33+66=",
    "5+31=", "2+81=", "71+9=", "23+60=", "29-13=",
    "3+62=", "1+13=", "4+68=", "41+0=", "33+31=",
    "57+35=", "21+21=", "20+29=", "6+3=", "92+2=",
    "32+3=", "9+27=", "91-4=", "71-67=", "56-2=",
    "91-11=", "53-3=", "15+74=", "22-1=", "30-1=",
    "74-4=", "4+81=", "60-27=", "13+14=", "14+40=",
    "94-44=", "32+18=", "36+34=", "77-63=", "77-41="
)

$t = $d.Tables(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output "Updated $i cells; date set."
